$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.99000000000062
$ws.Range("H2").Value = 0.00008380092456516941
$ws.Range("I2").Value = 0.00008380092456516941
$ws.Range("L2").Value = 56.042197172157941
$ws.Range("M2").Value = "[27.007339207122044, 85.07705513719384]"
$ws.Range("N2").Value = 0.00033034592812764002
$ws.Range("O2").Value = 0.00033034592812764002
$ws.Range("P2").Value = 1.742184514603349
$ws.Range("Q2").Value = "[1.0755001877154253, 2.408868841491273]"
$ws.Range("R2").Value = 0.0000038105689834644352
$ws.Range("S2").Value = 0.0000038105689834644352
$ws.Range("T2").Value = 68.921846522886014
$ws.Range("U2").Value = "[51.638357359024795, 86.20533568674723]"
$ws.Range("V2").Value = 0.00000000030998648092861458
$ws.Range("W2").Value = 0.00000000030998648092861458
$ws.Range("X2").Value = 18.783563563564009
$ws.Range("Y2").Value = 16.025865865866251
$ws.Range("Z2").Value = 21.54126126126177

# Row 3
$ws.Range("F3").Value = 25.99000000000062
$ws.Range("H3").Value = 0.00048352473308188593
$ws.Range("I3").Value = 0.00048352473308188593
$ws.Range("L3").Value = 48.170470000216973
$ws.Range("M3").Value = "[21.603971051204795, 74.73696894922915]"
$ws.Range("N3").Value = 0.00067583918406155519
$ws.Range("O3").Value = 0.00067583918406155519
$ws.Range("P3").Value = 1.528342372016656
$ws.Range("Q3").Value = "[0.7987632973091161, 2.257921446724196]"
$ws.Range("R3").Value = 0.00011713602124330261
$ws.Range("S3").Value = 0.00011713602124330261
$ws.Range("T3").Value = 64.194588472532956
$ws.Range("U3").Value = "[47.50718896491125, 80.88198798015466]"
$ws.Range("V3").Value = 0.00000000080399154001042916
$ws.Range("W3").Value = 0.00000000080399154001042916
$ws.Range("X3").Value = 19.66810810810858
$ws.Range("Y3").Value = 16.65025025025065
$ws.Range("Z3").Value = 22.68596596596651

# Row 4
$ws.Range("F4").Value = 25.99000000000062
$ws.Range("H4").Value = 0.000000080029015303573203
$ws.Range("I4").Value = 0.000000080029015303573203
$ws.Range("L4").Value = 71.835829470471126
$ws.Range("M4").Value = "[44.424468542037104, 99.24719039890515]"
$ws.Range("N4").Value = 0.000003622914280576595
$ws.Range("O4").Value = 0.000003622914280576595
$ws.Range("P4").Value = 0.94971069207619419
$ws.Range("Q4").Value = "[0.5597632555945777, 1.3396581285578106]"
$ws.Range("R4").Value = 0.00001260010382542198
$ws.Range("S4").Value = 0.00001260010382542198
$ws.Range("T4").Value = 65.683882175313713
$ws.Range("U4").Value = "[51.04377279704768, 80.32399155357975]"
$ws.Range("V4").Value = 0.000000000011256995335884311
$ws.Range("W4").Value = 0.000000000011256995335884311
$ws.Range("X4").Value = 22.061581581582111
$ws.Range("Y4").Value = 20.448588588589079
$ws.Range("Z4").Value = 23.67457457457515

# Row 5
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = 25.99000000000062
$ws.Range("H5").Value = 0.104922582719539
$ws.Range("I5").Value = 0.104922582719539
$ws.Range("L5").Value = 25.14265970155639
$ws.Range("M5").Value = "[-1.222148001512295, 51.50746740462507]"
$ws.Range("N5").Value = 0.061114391367280252
$ws.Range("O5").Value = 0.061114391367280252
$ws.Range("P5").Value = 0.15723686954903909
$ws.Range("Q5").Value = "[-1.446579199851156, 1.7610529389492342]"
$ws.Range("R5").Value = 0.84435594794701618
$ws.Range("S5").Value = 0.84435594794701618
$ws.Range("T5").Value = 75.39841911168493
$ws.Range("U5").Value = "[59.60006787250956, 91.1967703508603]"
$ws.Range("V5").Value = 0.000000000001767697099808174
$ws.Range("W5").Value = 0.000000000001767697099808174
$ws.Range("X5").Value = 25.339599599600209
$ws.Range("Y5").Value = 18.705515515515959
$ws.Range("Z5").Value = 31.973683683684449

# Row 6
$ws.Range("F6").Value = 23.09000000000017
$ws.Range("H6").Value = 0.000066862108230369977
$ws.Range("I6").Value = 0.000066862108230369977
$ws.Range("L6").Value = 58.647430220245823
$ws.Range("M6").Value = "[27.42452565529733, 89.87033478519432]"
$ws.Range("N6").Value = 0.00045467328570047982
$ws.Range("O6").Value = 0.00045467328570047982
$ws.Range("P6").Value = 0.42139481039142362
$ws.Range("Q6").Value = "[-0.1823947686768852, 1.0251843894597323]"
$ws.Range("R6").Value = 0.16668740500875101
$ws.Range("S6").Value = 0.16668740500875101
$ws.Range("T6").Value = 67.676747042973091
$ws.Range("U6").Value = "[50.438643309136225, 84.91485077680996]"
$ws.Range("V6").Value = 0.00000000047034287575797862
$ws.Range("W6").Value = 0.00000000047034287575797862
$ws.Range("X6").Value = 21.541421421421578
$ws.Range("Y6").Value = 19.322562562562709
$ws.Range("Z6").Value = 23.760280280280458

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 23.09000000000017
$ws.Range("H7").Value = 0.0054874799260604679
$ws.Range("I7").Value = 0.0054874799260604679
$ws.Range("L7").Value = 42.868395173024048
$ws.Range("M7").Value = "[7.857570905036461, 77.87921944101163]"
$ws.Range("N7").Value = 0.0175340269942863
$ws.Range("O7").Value = 0.0175340269942863
$ws.Range("P7").Value = 0.76102644861734792
$ws.Range("Q7").Value = "[0.08176317216550189, 1.440289725069194]"
$ws.Range("R7").Value = 0.02894115723803159
$ws.Range("S7").Value = 0.02894115723803159
$ws.Range("T7").Value = 62.056843715208231
$ws.Range("U7").Value = "[44.259831278680764, 79.8538561517357]"
$ws.Range("V7").Value = 0.0000000094195178323985829
$ws.Range("W7").Value = 0.0000000094195178323985829
$ws.Range("X7").Value = 20.293313313313458
$ws.Range("Y7").Value = 17.797097097097229
$ws.Range("Z7").Value = 22.789529529529691

# Row 8
$ws.Range("B8").Value = 1
$ws.Range("F8").Value = 23.09000000000017
$ws.Range("H8").Value = 0.002107766604203865
$ws.Range("I8").Value = 0.002107766604203865
$ws.Range("L8").Value = 43.372363361413633
$ws.Range("M8").Value = "[12.162801408445787, 74.58192531438146]"
$ws.Range("N8").Value = 0.007521921527160913
$ws.Range("O8").Value = 0.007521921527160913
$ws.Range("P8").Value = 0.87423699469265514
$ws.Range("Q8").Value = "[0.19497371824080734, 1.553500271144503]"
$ws.Range("R8").Value = 0.012814906377905279
$ws.Range("S8").Value = 0.012814906377905279
$ws.Range("T8").Value = 70.450409266621591
$ws.Range("U8").Value = "[54.13309120333, 86.76772732991319]"
$ws.Range("V8").Value = 0.000000000034227065626168951
$ws.Range("W8").Value = 0.000000000034227065626168951
$ws.Range("X8").Value = 19.877277277277429
$ws.Range("Y8").Value = 17.381061061061189
$ws.Range("Z8").Value = 22.373493493493658
